# Add "Wins"/"Losses"/"Ties" team-record columns (AD/AE/AF) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold/centered/bordered) from the last
# existing header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the same team record: 64 wins, 98 losses, 0 ties.
$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 98
    $ws.Cells.Item($r, 32).Value = 0
}
